$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1355
$ws.Range("I19").Value = 1355
$ws.Range("K19").Value = 1355
$ws.Range("M19").Value = -1180
$ws.Range("H58").Value = 415.2857
$ws.Range("I58").Value = 415.2857
$ws.Range("K58").Value = 1245.8571
$ws.Range("M58").Value = -1095.8571
$ws.Range("H86").Value = 361119360
$ws.Range("I86").Value = 285724900
$ws.Range("K86").Value = 285724900
$ws.Range("M86").Value = -285723777
$ws.Range("H89").Value = 361119360
$ws.Range("I89").Value = 285724900
$ws.Range("K89").Value = 1428624500
$ws.Range("M89").Value = -1428618884
$ws.Range("H106").Value = 13337133
$ws.Range("I106").Value = 16668916
$ws.Range("J106").Value = 10000
$ws.Range("K106").Value = 16668916
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = -16668285
$ws.Range("N106").Value = -11262
$ws.Range("H132").Value = 3104.889
$ws.Range("I132").Value = 2553.28
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 7659.84
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -5129.84
$ws.Range("N132").Value = -35060
$ws.Range("H137").Value = 29237.666
$ws.Range("I137").Value = 52749.562
$ws.Range("K137").Value = 158248.686
$ws.Range("M137").Value = -155698.686
$ws.Range("H138").Value = 4135.16
$ws.Range("I138").Value = 2682.1
$ws.Range("J138").Value = 4498.425
$ws.Range("K138").Value = 8046.299999999999
$ws.Range("L138").Value = 13495.275
$ws.Range("M138").Value = -2906.299999999999
$ws.Range("N138").Value = -23775.275
$ws.Range("H141").Value = 4432.485
$ws.Range("I141").Value = 4008.8064
$ws.Range("K141").Value = 12026.4192
$ws.Range("M141").Value = -6846.4192

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15406033
$ws.Range("I32").Value = 14939979
$ws.Range("J32").Value = 31251896
$ws.Range("K32").Value = 14939979
$ws.Range("L32").Value = 31251896
$ws.Range("M32").Value = -14939692
$ws.Range("N32").Value = -31252470
$ws.Range("H61").Value = 2776.56
$ws.Range("I61").Value = 2615.725
$ws.Range("J61").Value = 3419.9
$ws.Range("K61").Value = 2615.725
$ws.Range("L61").Value = 3419.9
$ws.Range("M61").Value = -2403.725
$ws.Range("N61").Value = -3843.9
$ws.Range("H110").Value = 1770.4839
$ws.Range("I110").Value = 1195.238
$ws.Range("J110").Value = 2978.5
$ws.Range("K110").Value = 1195.238
$ws.Range("L110").Value = 2978.5
$ws.Range("M110").Value = 849.7619999999999
$ws.Range("N110").Value = -7068.5
$ws.Range("H132").Value = 589886.8
$ws.Range("I132").Value = 715978.5
$ws.Range("J132").Value = 1459
$ws.Range("K132").Value = 2147935.5
$ws.Range("L132").Value = 4377
$ws.Range("M132").Value = -2145405.5
$ws.Range("N132").Value = -9437
$ws.Range("H136").Value = 2776.56
$ws.Range("I136").Value = 2615.725
$ws.Range("J136").Value = 3419.9
$ws.Range("K136").Value = 7847.174999999999
$ws.Range("L136").Value = 10259.7
$ws.Range("M136").Value = -5297.174999999999
$ws.Range("N136").Value = -15359.7

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3135.6667
$ws.Range("I99").Value = 2847.5557
$ws.Range("K99").Value = 2847.5557
$ws.Range("M99").Value = -1349.5557
$ws.Range("H134").Value = 1669527.2
$ws.Range("I134").Value = 1804302.6
$ws.Range("K134").Value = 5412907.800000001
$ws.Range("M134").Value = -5410372.800000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 139.9
$ws.Range("I7").Value = 152.66667
$ws.Range("K7").Value = 152.66667
$ws.Range("M7").Value = -39.66667000000001
$ws.Range("H31").Value = 4760.1313
$ws.Range("I31").Value = 1864.4445
$ws.Range("J31").Value = 5658.793
$ws.Range("K31").Value = 1864.4445
$ws.Range("L31").Value = 5658.793
$ws.Range("M31").Value = -1569.4445
$ws.Range("N31").Value = -6248.793
$ws.Range("H34").Value = 4760.1313
$ws.Range("I34").Value = 1864.4445
$ws.Range("J34").Value = 5658.793
$ws.Range("K34").Value = 1864.4445
$ws.Range("L34").Value = 5658.793
$ws.Range("M34").Value = -1662.4445
$ws.Range("N34").Value = -6062.793
$ws.Range("H86").Value = 9542.182000000001
$ws.Range("J86").Value = 9025.25
$ws.Range("L86").Value = 9025.25
$ws.Range("N86").Value = -11271.25
$ws.Range("H89").Value = 9542.182000000001
$ws.Range("J89").Value = 9025.25
$ws.Range("L89").Value = 45126.25
$ws.Range("N89").Value = -56358.25
$ws.Range("H134").Value = 2125.077
$ws.Range("I134").Value = 2140.1738
$ws.Range("K134").Value = 6420.5214
$ws.Range("M134").Value = -3885.5214

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 665
$ws.Range("I5").Value = 665
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1995
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1883
$ws.Range("N5").ClearContents()
$ws.Range("H52").Value = 4350.5
$ws.Range("J52").Value = 4350.5
$ws.Range("L52").Value = 13051.5
$ws.Range("N52").Value = -13583.5
$ws.Range("H68").Value = 3796.04
$ws.Range("I68").Value = 1974.2858
$ws.Range("J68").Value = 4504.5
$ws.Range("K68").Value = 5922.857400000001
$ws.Range("L68").Value = 13513.5
$ws.Range("M68").Value = -5111.857400000001
$ws.Range("N68").Value = -15135.5
$ws.Range("H71").Value = 3796.04
$ws.Range("I71").Value = 1974.2858
$ws.Range("J71").Value = 4504.5
$ws.Range("K71").Value = 17768.5722
$ws.Range("L71").Value = 40540.5
$ws.Range("M71").Value = -13712.5722
$ws.Range("N71").Value = -48652.5
$ws.Range("H103").Value = 562.8333
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4758
$ws.Range("H121").Value = 5082
$ws.Range("J121").Value = 6742.6665
$ws.Range("L121").Value = 20227.9995
$ws.Range("N121").Value = -22847.9995
$ws.Range("H131").Value = 1724.1714
$ws.Range("J131").Value = 1739.6111
$ws.Range("L131").Value = 5218.8333
$ws.Range("N131").Value = -15298.8333
$ws.Range("H135").Value = 665
$ws.Range("I135").Value = 665
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5985
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3450
$ws.Range("N135").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14927.786
$ws.Range("I70").Value = 33837.9
$ws.Range("J70").Value = 4422.1665
$ws.Range("K70").Value = 33837.9
$ws.Range("L70").Value = 4422.1665
$ws.Range("M70").Value = -33567.9
$ws.Range("N70").Value = -4962.1665
$ws.Range("H73").Value = 14927.786
$ws.Range("I73").Value = 33837.9
$ws.Range("J73").Value = 4422.1665
$ws.Range("K73").Value = 33837.9
$ws.Range("L73").Value = 4422.1665
$ws.Range("M73").Value = -32901.9
$ws.Range("N73").Value = -6294.1665
$ws.Range("H132").Value = 2560.8635
$ws.Range("I132").Value = 2414.5
$ws.Range("J132").Value = 2874.5
$ws.Range("K132").Value = 7243.5
$ws.Range("L132").Value = 8623.5
$ws.Range("M132").Value = -4713.5
$ws.Range("N132").Value = -13683.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9801.588
$ws.Range("J7").Value = 11938.875
$ws.Range("L7").Value = 11938.875
$ws.Range("N7").Value = -12162.875
$ws.Range("H16").Value = 3786.25
$ws.Range("I16").Value = 3786.25
$ws.Range("K16").Value = 3786.25
$ws.Range("M16").Value = -3616.25
$ws.Range("H40").Value = 66670452
$ws.Range("I40").Value = 83335560
$ws.Range("K40").Value = 83335560
$ws.Range("M40").Value = -83335424
$ws.Range("H93").Value = 58825468
$ws.Range("I93").Value = 83335140
$ws.Range("J93").Value = 2264
$ws.Range("K93").Value = 83335140
$ws.Range("L93").Value = 2264
$ws.Range("M93").Value = -83333892
$ws.Range("N93").Value = -4760
$ws.Range("H94").Value = 98950
$ws.Range("J94").Value = 98950
$ws.Range("L94").Value = 98950
$ws.Range("N94").Value = -100302
$ws.Range("H122").Value = 18708.092
$ws.Range("I122").Value = 20099.5
$ws.Range("K122").Value = 60298.5
$ws.Range("M122").Value = -57848.5
$ws.Range("H126").Value = 9801.588
$ws.Range("J126").Value = 11938.875
$ws.Range("L126").Value = 35816.625
$ws.Range("N126").Value = -40756.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1164.7273
$ws.Range("I100").Value = 1530.4286
$ws.Range("J100").Value = 524.75
$ws.Range("K100").Value = 3060.8572
$ws.Range("L100").Value = 1049.5
$ws.Range("M100").Value = -2519.8572
$ws.Range("N100").Value = -2131.5
$ws.Range("H126").Value = 4119.125
$ws.Range("I126").Value = 4119.125
$ws.Range("K126").Value = 12357.375
$ws.Range("M126").Value = -9887.375
$ws.Range("H132").Value = 28927.309
$ws.Range("I132").Value = 31841.771
$ws.Range("K132").Value = 95525.31299999999
$ws.Range("M132").Value = -92995.31299999999
$ws.Range("H136").Value = 23465.32
$ws.Range("I136").Value = 3199.3235
$ws.Range("K136").Value = 9597.970499999999
$ws.Range("M136").Value = -7047.970499999999
